# Update annual retention metrics (metricas_retencao_anual) to include the
# latest "visitar clientes" customer counts.
# For each affected cohort/period row, num_customers (column C) increases,
# and retention_rate (column E) = num_customers / cohort_size is recomputed.
# Row 37 (2024 cohort, period 0) also grows its cohort_size (column D)
# since num_customers == cohort_size at period 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: cohort 2021, period 5 -> num_customers 24 -> 25, cohort_size 2654
$ws.Cells.Item(22, 3).Value = 25
$ws.Cells.Item(22, 5).Value = 25 / 2654

# Row 27: cohort 2022, period 4 -> num_customers 36 -> 37, cohort_size 2252
$ws.Cells.Item(27, 3).Value = 37
$ws.Cells.Item(27, 5).Value = 37 / 2252

# Row 34: cohort 2023, period 2 -> num_customers 62 -> 68, cohort_size 2256
$ws.Cells.Item(34, 3).Value = 68
$ws.Cells.Item(34, 5).Value = 68 / 2256

# Row 36: cohort 2024, period 1 -> num_customers 99 -> 102, cohort_size 1930
$ws.Cells.Item(36, 3).Value = 102
$ws.Cells.Item(36, 5).Value = 102 / 1930

# Row 37: cohort 2025, period 0 -> num_customers and cohort_size both 604 -> 630
$ws.Cells.Item(37, 3).Value = 630
$ws.Cells.Item(37, 4).Value = 630
$ws.Cells.Item(37, 5).Value = 630 / 630
